# Update Excel files after daily scrape - 2025-12-16 03:32:44 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column width adjustments ----
$ws.Columns.Item(3).ColumnWidth = 56
$ws.Columns.Item(4).ColumnWidth = 39
$ws.Columns.Item(8).ColumnWidth = 29

# Column A (OPPORTUNITY ID) holds numeric-looking IDs that must stay text,
# exactly like the rest of the sheet (which uses inline strings throughout).
# Force text storage via a temporary "@" format, then drop the formatting
# again so the cells keep the sheet's default (unstyled) look.
$ws.Range("A2:A9").NumberFormat = "@"

# ---- Row 2 (existing row, values updated) ----
$ws.Cells.Item(2, 1).Value = "1326481"
$ws.Cells.Item(2, 2).Value = "https://aiesec.org/opportunity/global-talent/1326481"
$ws.Cells.Item(2, 3).Value = "Global Duty Billing Data Analytics Expert"
$ws.Cells.Item(2, 4).Value = "Maastricht, Netherlands"
$ws.Cells.Item(2, 5).Value = "Yes"
$ws.Cells.Item(2, 6).Value = "228 applicants"
$ws.Cells.Item(2, 7).Value = "6 - 18 Months"
$ws.Cells.Item(2, 8).Value = "DHL Group"

# ---- Row 3 (existing row, values updated) ----
$ws.Cells.Item(3, 1).Value = "1330623"
$ws.Cells.Item(3, 2).Value = "https://aiesec.org/opportunity/global-talent/1330623"
$ws.Cells.Item(3, 3).Value = "Intern – Strategy and Planning"
$ws.Cells.Item(3, 4).Value = "Karachi, Pakistan"
$ws.Cells.Item(3, 5).Value = "No"
$ws.Cells.Item(3, 6).Value = "0 applicants"
$ws.Cells.Item(3, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(3, 8).Value = "Keys Productions (Pvt) Ltd"

# ---- Row 4 (new) ----
$ws.Cells.Item(4, 1).Value = "1330617"
$ws.Cells.Item(4, 2).Value = "https://aiesec.org/opportunity/global-talent/1330617"
$ws.Cells.Item(4, 3).Value = "Content Marketing Intern"
$ws.Cells.Item(4, 4).Value = "Visakhapatnam, Andhra Pradesh, India"
$ws.Cells.Item(4, 5).Value = "No"
$ws.Cells.Item(4, 6).Value = "0 applicants"
$ws.Cells.Item(4, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(4, 8).Value = "KEN FOUNDATION SOCIETY"

# ---- Row 5 (new) ----
$ws.Cells.Item(5, 1).Value = "1330616"
$ws.Cells.Item(5, 2).Value = "https://aiesec.org/opportunity/global-talent/1330616"
$ws.Cells.Item(5, 3).Value = "Dental Assistant"
$ws.Cells.Item(5, 4).Value = "Visakhapatnam, Andhra Pradesh, India"
$ws.Cells.Item(5, 5).Value = "No"
$ws.Cells.Item(5, 6).Value = "0 applicants"
$ws.Cells.Item(5, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(5, 8).Value = "Kaushal's dental care"

# ---- Row 6 (new) ----
$ws.Cells.Item(6, 1).Value = "1329349"
$ws.Cells.Item(6, 2).Value = "https://aiesec.org/opportunity/global-talent/1329349"
$ws.Cells.Item(6, 3).Value = "Copywriter"
$ws.Cells.Item(6, 4).Value = "Yerevan, Armenia"
$ws.Cells.Item(6, 5).Value = "No"
$ws.Cells.Item(6, 6).Value = "80 applicants"
$ws.Cells.Item(6, 7).Value = "6 - 18 Months"
$ws.Cells.Item(6, 8).Value = "TCF Armenia"

# ---- Row 7 (new) ----
$ws.Cells.Item(7, 1).Value = "1328030"
$ws.Cells.Item(7, 2).Value = "https://aiesec.org/opportunity/global-talent/1328030"
$ws.Cells.Item(7, 3).Value = "[Remote] Software Application Support and Development"
$ws.Cells.Item(7, 4).Value = "No location available"
$ws.Cells.Item(7, 5).Value = "No"
$ws.Cells.Item(7, 6).Value = "23 applicants"
$ws.Cells.Item(7, 7).Value = "Remote"
$ws.Cells.Item(7, 8).Value = "dJava Factory Sdn Bhd"

# ---- Row 8 (new) ----
$ws.Cells.Item(8, 1).Value = "1328026"
$ws.Cells.Item(8, 2).Value = "https://aiesec.org/opportunity/global-talent/1328026"
$ws.Cells.Item(8, 3).Value = "[Remote] Software Application Support and Development"
$ws.Cells.Item(8, 4).Value = "No location available"
$ws.Cells.Item(8, 5).Value = "No"
$ws.Cells.Item(8, 6).Value = "35 applicants"
$ws.Cells.Item(8, 7).Value = "Remote"
$ws.Cells.Item(8, 8).Value = "dJava Factory Sdn Bhd"

# ---- Row 9 (new) ----
$ws.Cells.Item(9, 1).Value = "1309734"
$ws.Cells.Item(9, 2).Value = "https://aiesec.org/opportunity/global-talent/1309734"
$ws.Cells.Item(9, 3).Value = "Marketing Intern"
$ws.Cells.Item(9, 4).Value = "Lahore, Punjab, Pakistan"
$ws.Cells.Item(9, 5).Value = "No"
$ws.Cells.Item(9, 6).Value = "12 applicants"
$ws.Cells.Item(9, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(9, 8).Value = "Devsinc."

# Drop the temporary text-number-format now that the IDs are safely stored
# as text, restoring the plain/default cell style used by the rest of the sheet.
$ws.Range("A2:A9").ClearFormats()
